$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Ana Milic"
$ws.Range("B8").Value = "067/777-888"
$ws.Range("C8").Value = "TV"
$ws.Range("D8").Value = "frižider komb"
$ws.Range("E8").Value = "Samsung"
$ws.Range("F8").Value = "RF-850"
$ws.Range("G8").Value = "'22222"
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Value = "ne hladi zamrzivač"
